# Update the "want to go" counts (column F) on several sheets to reflect
# the latest generated data (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 215
$ws.Range("F4").Value = 409
$ws.Range("F6").Value = 5596
$ws.Range("F8").Value = 711
$ws.Range("F15").Value = 24
$ws.Range("F17").Value = 1885
$ws.Range("F18").Value = 1482
$ws.Range("F24").Value = 163
$ws.Range("F28").Value = 3025
$ws.Range("F33").Value = 40
$ws.Range("F34").Value = 416
$ws.Range("F40").Value = 744
$ws.Range("F41").Value = 95

# Sheet "演出" (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 207
$ws.Range("F6").Value = 145

# Sheet "全部类型" (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 215
$ws.Range("F7").Value = 5596
$ws.Range("F9").Value = 711
$ws.Range("F11").Value = 207
$ws.Range("F15").Value = 145
$ws.Range("F20").Value = 24
$ws.Range("F23").Value = 1885
$ws.Range("F24").Value = 1482
$ws.Range("F30").Value = 163
$ws.Range("F32").Value = 3025
$ws.Range("F37").Value = 40
$ws.Range("F38").Value = 416
$ws.Range("F43").Value = 744
$ws.Range("F44").Value = 95
